$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.429.32"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "1.655.25"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.00%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "213.71"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -1.80%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "24.16"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "1.888.80"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "1.650.98"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("E15").Value = "  +2.49%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "65.88"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "27.438.73"
$ws.Range("E17").Value = "  -2.11%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "232.39"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -7.48%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("E24").Value = "  -1.90%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "146.73"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "7.22"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "15.95"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -1.32%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.20"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "1.463.25"
$ws.Range("E33").Value = "  +2.55%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "3.11"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("E36").Value = "  +0.05%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.911"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +0.05%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "5.47"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "65.52"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "1.797.33"
$ws.Range("E45").Value = "  -2.01%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.781"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -2.15%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.74"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "88.42"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -0.37%  "
